# Insert a new weekly price record for "Ajo" (Primera/Chino) at Macroferia
# Regional de Talca, pushing the existing historical rows (497:552) down by
# one row to 498:553, then populate the freshly inserted row 497 with the
# new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 497:552 down to 498:553, leaving a blank row 497 (row height /
# per-cell formatting - e.g. the date style on column D - is carried along
# by Insert, matching Excel's native "insert copies format from above").
$ws.Rows(497).Insert()

# Fill the newly inserted row with the new reporting week's data.
$ws.Range("A497").Value = 5
$ws.Range("B497").Value = "Macroferia Regional de Talca"
$ws.Range("C497").Value = "Maule"
$ws.Range("D497").Value = 45212
$ws.Range("E497").Value = 7
$ws.Range("F497").Value = 100112003
$ws.Range("G497").Value = "Ajo"
$ws.Range("H497").Value = "Chino"
$ws.Range("I497").Value = "Primera"
$ws.Range("J497").Value = 150
$ws.Range("K497").Value = 21000
$ws.Range("L497").Value = 21000
$ws.Range("M497").Value = 21000
$ws.Range("N497").Value = "$/malla 10 kilos"
$ws.Range("O497").Value = "China"
$ws.Range("P497").Value = 2100
$ws.Range("Q497").Value = 10
$ws.Range("R497").Value = "Hortaliza"
